{"js": "// The document contains two places where a transcribed \"<id>...</id>\" tag\n// was split across three runs: a formatted \"<id>\" run, a plain run holding\n// the old record id, and a formatted \"</id>\" run. This edit collapses each\n// trio into a single run carrying the complete, updated id so the id text\n// reads as one contiguous, consistently-formatted run.\nconst replacements = [\n  { oldText: \"<id>p086v_a1</id>\", newText: \"<id>p086v_1</id>\" },\n  { oldText: \"<id>p087r_a1</id>\", newText: \"<id>p087r_1</id>\" }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  // Replacing the whole matched range merges the three original runs into\n  // one run that keeps the formatting of the first (leading \"<id>\") run.\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document contains two places where a transcribed \"<id>...</id>\" tag\n# was split across three runs: a formatted \"<id>\" run, a plain run holding\n# the old record id, and a formatted \"</id>\" run. This edit collapses each\n# trio into a single run carrying the complete, updated id so the id text\n# reads as one contiguous, consistently-formatted run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"<id>p086v_a1</id>\"; New = \"<id>p086v_1</id>\" },\n    @{ Old = \"<id>p087r_a1</id>\"; New = \"<id>p087r_1</id>\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = \"wdFindContinue\"\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # Replacing the whole matched range merges the three original runs into\n    # one run that keeps the formatting of the first (leading \"<id>\") run.\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, \"wdFindContinue\", $false, $pair.New, \"wdReplaceAll\")\n}\n"}
